$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$c = $ws.Range("D2")
$c.NumberFormat = "@"
$c.Value = '45.839.48'
$c.Style = "Normal"

$ws.Range("E2").Value = '  -2.12%  '

$c = $ws.Range("D3")
$c.NumberFormat = "@"
$c.Value = '2.370.83'
$c.Style = "Normal"

$ws.Range("E3").Value = '  +2.73%  '

$ws.Range("E4").Value = '  -0.02%  '

$c = $ws.Range("D5")
$c.NumberFormat = "@"
$c.Value = '300.04'
$c.Style = "Normal"

$ws.Range("E5").Value = '  -1.41%  '

$c = $ws.Range("D6")
$c.NumberFormat = "@"
$c.Value = '98.36'
$c.Style = "Normal"

$ws.Range("E6").Value = '  -3.28%  '

$c = $ws.Range("D7")
$c.NumberFormat = "@"
$c.Value = '0.566'
$c.Style = "Normal"

$ws.Range("E7").Value = '  -1.07%  '

$ws.Range("E8").Value = '  +0.01%  '

$c = $ws.Range("D10")
$c.NumberFormat = "@"
$c.Value = '34.23'
$c.Style = "Normal"

$ws.Range("E10").Value = '  -7.51%  '

$c = $ws.Range("D11")
$c.NumberFormat = "@"
$c.Value = '0.0786'
$c.Style = "Normal"

$ws.Range("E11").Value = '  -2.21%  '

$c = $ws.Range("D12")
$c.NumberFormat = "@"
$c.Value = '7.09'
$c.Style = "Normal"

$ws.Range("E12").Value = '  -4.56%  '

$ws.Range("E13").Value = '  -0.40%  '

$c = $ws.Range("D14")
$c.NumberFormat = "@"
$c.Value = '2.734.25'
$c.Style = "Normal"

$ws.Range("E14").Value = '  +2.83%  '

$c = $ws.Range("D15")
$c.NumberFormat = "@"
$c.Value = '2.382.50'
$c.Style = "Normal"

$ws.Range("E15").Value = '  +3.20%  '

$c = $ws.Range("D16")
$c.NumberFormat = "@"
$c.Value = '0.820'
$c.Style = "Normal"

$ws.Range("E16").Value = '  -0.13%  '

$c = $ws.Range("D17")
$c.NumberFormat = "@"
$c.Value = '13.69'
$c.Style = "Normal"

$ws.Range("E17").Value = '  -2.47%  '

$c = $ws.Range("D18")
$c.NumberFormat = "@"
$c.Value = '45.780.29'
$c.Style = "Normal"

$ws.Range("E18").Value = '  -2.23%  '

$c = $ws.Range("D19")
$c.NumberFormat = "@"
$c.Value = '12.70'
$c.Style = "Normal"

$ws.Range("E19").Value = '  -8.44%  '

$c = $ws.Range("D20")
$c.NumberFormat = "@"
$c.Value = '0.0₃0946'
$c.Style = "Normal"

$ws.Range("E20").Value = '  -0.44%  '

$c = $ws.Range("D21")
$c.NumberFormat = "@"
$c.Value = '6.02'
$c.Style = "Normal"

$ws.Range("E21").Value = '  -1.57%  '

$c = $ws.Range("D22")
$c.NumberFormat = "@"
$c.Value = '66.75'
$c.Style = "Normal"

$ws.Range("E22").Value = '  -0.28%  '

$c = $ws.Range("D23")
$c.NumberFormat = "@"
$c.Value = '243.56'
$c.Style = "Normal"

$ws.Range("E23").Value = '  -2.46%  '

$ws.Range("E24").Value = '  -5.46%  '

$ws.Range("E25").Value = '  -0.06%  '

$c = $ws.Range("D26")
$c.NumberFormat = "@"
$c.Value = '1.91'
$c.Style = "Normal"

$ws.Range("E26").Value = '  -2.63%  '

$c = $ws.Range("D27")
$c.NumberFormat = "@"
$c.Value = '38.80'
$c.Style = "Normal"

$ws.Range("E27").Value = '  -11.78%  '

$ws.Range("E28").Value = '  -3.55%  '

$ws.Range("E29").Value = '  -2.93%  '

$c = $ws.Range("D30")
$c.NumberFormat = "@"
$c.Value = '20.92'
$c.Style = "Normal"

$ws.Range("E30").Value = '  +3.35%  '

$c = $ws.Range("D31")
$c.NumberFormat = "@"
$c.Value = '3.72'
$c.Style = "Normal"

$ws.Range("E31").Value = '  +15.84%  '

$ws.Range("E32").Value = '  -2.55%  '

$ws.Range("E33").Value = '  -5.14%  '

$c = $ws.Range("D34")
$c.NumberFormat = "@"
$c.Value = '147.30'
$c.Style = "Normal"

$ws.Range("E34").Value = '  -0.44%  '

$ws.Range("E35").Value = '  -5.32%  '

$ws.Range("E36").Value = '  -1.23%  '

$c = $ws.Range("D37")
$c.NumberFormat = "@"
$c.Value = '1.90'
$c.Style = "Normal"

$ws.Range("E37").Value = '  +4.43%  '

$c = $ws.Range("D38")
$c.NumberFormat = "@"
$c.Value = '0.116'
$c.Style = "Normal"

$ws.Range("E38").Value = '  -3.20%  '

$c = $ws.Range("D39")
$c.NumberFormat = "@"
$c.Value = '14.93'
$c.Style = "Normal"

$ws.Range("E39").Value = '  -8.20%  '

$c = $ws.Range("D40")
$c.NumberFormat = "@"
$c.Value = '3.83'
$c.Style = "Normal"

$ws.Range("E40").Value = '  -4.99%  '

$ws.Range("E41").Value = '  -3.26%  '

$ws.Range("E42").Value = '  -8.18%  '

$c = $ws.Range("D43")
$c.NumberFormat = "@"
$c.Value = '1.941.98'
$c.Style = "Normal"

$ws.Range("E43").Value = '  +4.46%  '

$c = $ws.Range("D44")
$c.NumberFormat = "@"
$c.Value = '0.999'
$c.Style = "Normal"

$ws.Range("E44").Value = '  +0.05%  '

$c = $ws.Range("D45")
$c.NumberFormat = "@"
$c.Value = '95.12'
$c.Style = "Normal"

$ws.Range("E45").Value = '  +7.02%  '

$ws.Range("E46").Value = '  -10.73%  '

$ws.Range("E47").Value = '  +5.58%  '

$ws.Range("E48").Value = '  -5.84%  '

$c = $ws.Range("D49")
$c.NumberFormat = "@"
$c.Value = '98.89'
$c.Style = "Normal"

$ws.Range("E49").Value = '  +1.59%  '

$c = $ws.Range("D50")
$c.NumberFormat = "@"
$c.Value = '2.604.66'
$c.Style = "Normal"

$ws.Range("E50").Value = '  +2.74%  '

$c = $ws.Range("D51")
$c.NumberFormat = "@"
$c.Value = '68.54'
$c.Style = "Normal"

$ws.Range("E51").Value = '  -8.38%  '
